{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Target: the paragraph that discusses the \"SINGLE_LINKED\" vs \"ARRAY_LIST\"\n// timing comparison. The original text:\n//   \u201cSINGLE_LINKED\u201d se demor\u00f3 20 segundos m\u00e1s; no obstante, consideramos que\n//   esta diferencia de tiempo no es significativa en este caso.\n// is replaced by an expanded conclusion paragraph (split across several\n// runs, all sharing the same \"Dax-Regular\" / es-419 run formatting that the\n// surrounding text already uses).\n\nconst body = context.document.body;\n\n// The exact original run text (3 runs concatenated) that is being replaced.\nconst originalText =\n  \"\\u201cSINGLE_LINKED\\u201d se demor\\u00f3 20 segundos m\\u00e1s; no obstante, \" +\n  \"consideramos que esta diferencia de tiempo no es significativa en este caso.\";\n\nconst results = body.search(originalText, { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly 1 match for the target sentence, found \" + results.items.length\n  );\n}\n\nconst targetRange = results.items[0];\n\n// New run texts (in order) - matches the diff's resulting <w:r> sequence.\n// Entries that need leading/trailing spaces preserved use xml:space=\"preserve\".\nconst newRuns = [\n  { text: \"\\u201cSINGLE_LINKED\\u201d \", preserve: true },\n  { text: \"se demor\\u00f3 20 segundos m\\u00e1s\", preserve: false },\n  {\n    text:\n      \". Teniendo en cuenta estos resultados, se puede concluir que utilizar un \" +\n      \"\\u201cARRAY_LIST\\u201d\",\n    preserve: false,\n  },\n  {\n    text:\n      \"puede llegar a ser mucho m\\u00e1s eficiente (en casos similares al del \" +\n      \"laboratorio) cuando se trata un gran n\\u00famero de datos, es por esto que \" +\n      \"entre m\\u00e1s grande sea la cantidad de datos, mayor ser\\u00e1 la diferencia \" +\n      \"de tiempo entre un\",\n    preserve: false,\n  },\n  { text: \"a lista de tipo\", preserve: false },\n  { text: \" \\u201cARRAY_LIST\\u201d\", preserve: true },\n  { text: \" y una de \\u201cSINGLE_LINKED\\u201d.\", preserve: true },\n];\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\nconst runProps =\n  '<w:rPr><w:rFonts w:ascii=\"Dax-Regular\" w:hAnsi=\"Dax-Regular\"/>' +\n  '<w:lang w:val=\"es-419\"/></w:rPr>';\n\nconst runsXml = newRuns\n  .map((run) => {\n    const spaceAttr = run.preserve ? ' xml:space=\"preserve\"' : \"\";\n    return (\n      \"<w:r>\" +\n      runProps +\n      \"<w:t\" +\n      spaceAttr +\n      \">\" +\n      escapeXml(run.text) +\n      \"</w:t></w:r>\"\n    );\n  })\n  .join(\"\");\n\nconst ooxmlPackage =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body><w:p>\" +\n  runsXml +\n  \"</w:p></w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\ntargetRange.insertOoxml(ooxmlPackage, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Target: the paragraph that discusses the \"SINGLE_LINKED\" vs \"ARRAY_LIST\"\n# timing comparison. The original text:\n#   \u201cSINGLE_LINKED\u201d se demor\u00f3 20 segundos m\u00e1s; no obstante, consideramos que\n#   esta diferencia de tiempo no es significativa en este caso.\n# is replaced by an expanded conclusion paragraph (split across several\n# runs, all sharing the same \"Dax-Regular\" / es-419 run formatting that the\n# surrounding text already uses).\n\n$d = $word.ActiveDocument\n\n# Exact original text (spanning the 3 original runs) that we are replacing.\n$originalText = \"\u201cSINGLE_LINKED\u201d se demor\u00f3 20 segundos m\u00e1s; no obstante, consideramos que esta diferencia de tiempo no es significativa en este caso.\"\n\n# Locate it (use a duplicate range so we don't disturb $d.Content itself).\n$find = $d.Content.Duplicate\n$find.Find.MatchCase = $true\n$find.Find.Text = $originalText\n$found = $find.Find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the target sentence to replace.\"\n}\n\n# New run texts, in order, mirroring the diff's resulting <w:r> sequence.\n# Entries with leading/trailing spaces need xml:space=\"preserve\".\n$newRuns = @(\n    @{ Text = \"\u201cSINGLE_LINKED\u201d \"; Preserve = $true },\n    @{ Text = \"se demor\u00f3 20 segundos m\u00e1s\"; Preserve = $false },\n    @{ Text = \". Teniendo en cuenta estos resultados, se puede concluir que utilizar un \u201cARRAY_LIST\u201d\"; Preserve = $false },\n    @{ Text = \"puede llegar a ser mucho m\u00e1s eficiente (en casos similares al del laboratorio) cuando se trata un gran n\u00famero de datos, es por esto que entre m\u00e1s grande sea la cantidad de datos, mayor ser\u00e1 la diferencia de tiempo entre un\"; Preserve = $false },\n    @{ Text = \"a lista de tipo\"; Preserve = $false },\n    @{ Text = \" \u201cARRAY_LIST\u201d\"; Preserve = $true },\n    @{ Text = \" y una de \u201cSINGLE_LINKED\u201d.\"; Preserve = $true }\n)\n\nfunction Escape-Xml([string]$s) {\n    return $s.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\")\n}\n\n$runProps = '<w:rPr><w:rFonts w:ascii=\"Dax-Regular\" w:hAnsi=\"Dax-Regular\"/><w:lang w:val=\"es-419\"/></w:rPr>'\n\n$runsXml = \"\"\nforeach ($run in $newRuns) {\n    $spaceAttr = \"\"\n    if ($run.Preserve) {\n        $spaceAttr = ' xml:space=\"preserve\"'\n    }\n    $runsXml += \"<w:r>\" + $runProps + \"<w:t\" + $spaceAttr + \">\" + (Escape-Xml $run.Text) + \"</w:t></w:r>\"\n}\n\n$ooxmlPackage = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' + $runsXml + '</w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n# Build a fresh Range over exactly the matched text (by character offsets)\n# and replace its content in one shot via InsertXML, so the paragraph stays\n# intact (no extra paragraph break is introduced).\n$targetRange = $d.Range($find.Start, $find.End)\n$targetRange.InsertXML($ooxmlPackage)\n"}
